# Recorte de carga en parametros
$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet, so it lands at the end
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Recorte Carga"

# Fill in the headers (plain text, no special style)
$newSheet.Range("A1").Value = "Bus"
$newSheet.Range("B1").Value = "Recorte Max"

# Make the new sheet the active/selected tab
$newSheet.Activate()
